# Refresh the cryptos price list (Price + Volume(1h) columns) with the
# latest values from the GitHub Actions scrape run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds values such as "3.105.73" or "0.999" which look
# numeric but must stay plain text (that is how the sheet already stores
# them). Writing straight into Range.Value lets Excel "helpfully" coerce
# strings that parse as a number (e.g. "1.00" or "560.34") into real
# numbers, which would change both the stored type and the display value.
# To keep the cells as genuine text -- without permanently changing the
# number format of the target cells themselves -- stage each value in a
# scratch cell formatted as Text, then copy/PasteSpecial just the value
# into the destination; PasteSpecial(xlPasteValues) carries over the
# scratch cells text value but leaves the destinations own formatting
# untouched.
$xlPasteValues = -4163
$scratch = $ws.Range("Z100")
$scratch.NumberFormat = "@"

function Set-PriceText($cellRef, $text) {
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial($xlPasteValues)
}

Set-PriceText "D2" "65.009.29"
$ws.Range("E2").Value = "  +3.80%  "
Set-PriceText "D3" "3.105.16"
$ws.Range("E3").Value = "  +1.88%  "
Set-PriceText "D4" "1.00"
$ws.Range("E4").Value = "  +0.06%  "
Set-PriceText "D5" "560.34"
$ws.Range("E5").Value = "  +2.52%  "
Set-PriceText "D6" "144.97"
$ws.Range("E6").Value = "  +6.81%  "
$ws.Range("E7").Value = "  -0.06%  "
Set-PriceText "D8" "3.106.62"
$ws.Range("E8").Value = "  +2.19%  "
Set-PriceText "D9" "0.502"
$ws.Range("E9").Value = "  +1.22%  "
Set-PriceText "D10" "7.19"
$ws.Range("E10").Value = "  +17.18%  "
Set-PriceText "D11" "0.153"
$ws.Range("E11").Value = "  +2.84%  "
Set-PriceText "D12" "0.468"
$ws.Range("E12").Value = "  +3.81%  "
$ws.Range("E13").Value = "  +4.94%  "
Set-PriceText "D14" "35.49"
$ws.Range("E14").Value = "  +1.30%  "
Set-PriceText "D15" "3.607.39"
$ws.Range("E15").Value = "  +1.68%  "
Set-PriceText "D16" "65.170.16"
$ws.Range("E16").Value = "  +3.89%  "
Set-PriceText "D17" "3.104.51"
$ws.Range("E17").Value = "  +1.66%  "
Set-PriceText "D19" "6.85"
$ws.Range("E19").Value = "  +2.32%  "
Set-PriceText "D20" "482.39"
$ws.Range("E20").Value = "  -0.46%  "
Set-PriceText "D21" "13.83"
$ws.Range("E21").Value = "  +3.62%  "
Set-PriceText "D22" "0.681"
$ws.Range("E22").Value = "  +0.18%  "
Set-PriceText "D23" "7.53"
$ws.Range("E23").Value = "  +6.25%  "
Set-PriceText "D24" "13.52"
$ws.Range("E24").Value = "  +11.28%  "
Set-PriceText "D25" "81.31"
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("E26").Value = "  +0.11%  "
Set-PriceText "D27" "2.79"
$ws.Range("E27").Value = "  +2.25%  "
Set-PriceText "D28" "8.21"
$ws.Range("E28").Value = "  +4.76%  "
Set-PriceText "D29" "2.06"
$ws.Range("E29").Value = "  +5.70%  "
$ws.Range("E30").Value = "  -0.27%  "
Set-PriceText "D31" "26.16"
$ws.Range("E31").Value = "  +1.02%  "
$ws.Range("E32").Value = "  +1.66%  "
Set-PriceText "D33" "2.49"
$ws.Range("E33").Value = "  +5.04%  "
Set-PriceText "D34" "5.68"
$ws.Range("E34").Value = "  -0.78%  "
Set-PriceText "D35" "6.24"
$ws.Range("E35").Value = "  +5.41%  "
Set-PriceText "D36" "55.23"
$ws.Range("E36").Value = "  +0.15%  "
Set-PriceText "D37" "471.39"
$ws.Range("E37").Value = "  +2.72%  "
Set-PriceText "D38" "0.0411"
$ws.Range("E38").Value = "  +6.08%  "
Set-PriceText "D39" "0.0832"
$ws.Range("E39").Value = "  +3.16%  "
Set-PriceText "D40" "2.92"
$ws.Range("E40").Value = "  +19.28%  "
Set-PriceText "D41" "3.013.06"
$ws.Range("E41").Value = "  -5.41%  "
Set-PriceText "D42" "8.28"
$ws.Range("E42").Value = "  +1.50%  "
Set-PriceText "D43" "0.116"
$ws.Range("E43").Value = "  -2.00%  "
Set-PriceText "D44" "28.21"
$ws.Range("E44").Value = "  +6.23%  "
$ws.Range("E45").Value = "  +5.19%  "
$ws.Range("E46").Value = "  +0.06%  "
Set-PriceText "D47" "2.13"
$ws.Range("E47").Value = "  +7.37%  "
$ws.Range("E48").Value = "  +2.80%  "
Set-PriceText "D49" "0.0₃0522"
$ws.Range("E49").Value = "  +5.74%  "
Set-PriceText "D50" "116.58"
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("E51").Value = "  +2.08%  "

# Clean up the scratch cell so it leaves no trace in the saved sheet.
$scratch.Clear()
